$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns keep their text formatting
# so numeric-looking strings (e.g. "37.50", "1.18%") are preserved verbatim
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '310.98'
$ws.Range("E2").Value = '1.18%'
$ws.Range("D3").Value = '37.50'
$ws.Range("E3").Value = '0.09%'
$ws.Range("D4").Value = '5.107'
$ws.Range("E4").Value = '0.12%'
$ws.Range("D5").Value = '0.07853'
$ws.Range("E5").Value = '-0.24%'
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").Value = '4.392'
$ws.Range("E6").Value = '1.07%'
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").Value = '1.903'
$ws.Range("E7").Value = '-3.99%'
$ws.Range("B8").Value = 'KuCoinToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D8").Value = '8.216'
$ws.Range("E8").Value = '-0.12%'
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").Value = '2.724'
$ws.Range("E9").Value = '-13.02%'
$ws.Range("B10").Value = 'MXToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D10").Value = '0.9268'
$ws.Range("E10").Value = '0.09%'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = '0.1198'
$ws.Range("E11").Value = '-7.05%'
$ws.Range("B12").Value = 'WazirX'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D12").Value = '0.1903'
$ws.Range("E12").Value = '0.45%'
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D13").Value = '0.09396'
$ws.Range("E13").Value = '5.67%'
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").Value = '0.03416'
$ws.Range("E14").Value = '-0.26%'
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").Value = '0.09618'
$ws.Range("E15").Value = '-1.38%'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Value = '0.001370'
$ws.Range("E16").Value = '-1.79%'
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").Value = '0.005855'
$ws.Range("E17").Value = '-2.12%'
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").Value = '3.537'
$ws.Range("E18").Value = '-1.23%'
$ws.Range("E19").Value = '-0.28%'
$ws.Range("D20").Value = '5.255'
$ws.Range("E20").Value = '4.91%'
$ws.Range("D21").Value = '0.1274'
$ws.Range("E21").Value = '-0.70%'
$ws.Range("D22").Value = '0.2583'
$ws.Range("E22").Value = '3.51%'
$ws.Range("D23").Value = '0.02099'
$ws.Range("E23").Value = '179.64%'
$ws.Range("D24").Value = '0.04345'
$ws.Range("E24").Value = '0.52%'
$ws.Range("D25").Value = '0.001194'
$ws.Range("E25").Value = '-2.30%'
$ws.Range("D26").Value = '0.004272'
$ws.Range("E26").Value = '-7.09%'
$ws.Range("D27").Value = '0.0001297'
$ws.Range("E27").Value = '-63.91%'
$ws.Range("D39").Value = '0.02093'
$ws.Range("E39").Value = '-8.65%'
$ws.Range("D40").Value = '0.05047'
$ws.Range("E40").Value = '0.72%'
$ws.Range("D41").Value = '0.007628'
$ws.Range("E41").Value = '1.59%'
$ws.Range("D42").Value = '0.009107'
$ws.Range("E42").Value = '-7.79%'
$ws.Range("D43").Value = '0.1350'
$ws.Range("E43").Value = '-0.20%'
$ws.Range("D44").Value = '0.002000'
$ws.Range("E44").Value = '-4.47%'
$ws.Range("D45").Value = '0.008573'
$ws.Range("E45").Value = '6.90%'
$ws.Range("D46").Value = '0.00006690'
$ws.Range("E46").Value = '2.58%'
$ws.Range("D47").Value = '0.00000000748'
$ws.Range("E47").Value = '-0.45%'
$ws.Range("D48").Value = '0.001197'
$ws.Range("E48").Value = '-0.42%'
$ws.Range("D49").Value = '0.002896'
$ws.Range("E49").Value = '-3.57%'
$ws.Range("D50").Value = '0.00002095'
$ws.Range("E50").Value = '-0.45%'
$ws.Range("D51").Value = '0.0001995'
$ws.Range("E51").Value = '-0.45%'
